{"js": "// Update in-text citation placeholder codes throughout the document body.\n// Each old \"Ref-XXXXXXX\" marker is replaced with its corresponding new\n// marker, per the dataset/citation-check update described in the commit.\nconst replacements = [\n  [\"Ref-JHG7Y6\", \"Ref-u164053\"],\n  [\"Ref-K89UIO\", \"Ref-u164053\"],\n  [\"Ref-DJ49F2\", \"Ref-s937590\"],\n  [\"Ref-J7X2B9\", \"Ref-u865466\"],\n  [\"Ref-J49F2K\", \"Ref-u192333\"],\n  [\"Ref-DJ49KL\", \"Ref-f811037\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Update in-text citation placeholder codes throughout the document body.\n# Each old \"Ref-XXXXXXX\" marker is replaced with its corresponding new\n# marker, per the dataset/citation-check update described in the commit.\n\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @{ Old = \"Ref-JHG7Y6\"; New = \"Ref-u164053\" },\n    @{ Old = \"Ref-K89UIO\"; New = \"Ref-u164053\" },\n    @{ Old = \"Ref-DJ49F2\"; New = \"Ref-s937590\" },\n    @{ Old = \"Ref-J7X2B9\"; New = \"Ref-u865466\" },\n    @{ Old = \"Ref-J49F2K\"; New = \"Ref-u192333\" },\n    @{ Old = \"Ref-DJ49KL\"; New = \"Ref-f811037\" }\n)\n\nforeach ($pair in $replacements) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $pair.Old\n    $find.Replacement.Text = $pair.New\n    $find.Forward = $true\n    $find.Wrap = 1\n    $find.Execute([ref]$pair.Old, $false, $false, $false, $false, $false, $true, 1, $false, [ref]$pair.New, 2)\n}\n"}
